$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '62.897.44'
$ws.Cells.Item(2, 5).Value = '  +6.22%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.446.35'
$ws.Cells.Item(3, 5).Value = '  +3.73%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.13%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '617.26'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +11.18%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '145.86'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +6.79%  '

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -0.25%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +2.04%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.445.16'
$ws.Cells.Item(9, 5).Value = '  +4.02%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +6.28%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.31%  '

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '5.24'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +3.94%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +5.21%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +6.47%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +9.46%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.973.60'
$ws.Cells.Item(16, 5).Value = '  +6.90%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '62.643.56'
$ws.Cells.Item(17, 5).Value = '  +5.61%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.445.17'
$ws.Cells.Item(18, 5).Value = '  +3.90%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.61%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +5.45%  '

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '326.68'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +2.02%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +3.34%  '

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '2.04'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +14.76%  '

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.15%  '

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '65.79'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +3.01%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '1.17'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +16.89%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Bittensor'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '616.07'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +12.39%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Aptos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '8.36'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +4.24%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).Value = '0.0₃0982'
$ws.Cells.Item(29, 5).Value = '  +8.35%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'WrappedeETH'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(30, 4).Value = '2.565.06'
$ws.Cells.Item(30, 5).Value = '  +3.81%  '

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '8.14'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +3.11%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +9.84%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '1.86'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +6.01%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'Kaspa'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '0.137'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +6.43%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +5.79%  '

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -0.23%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +5.63%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +2.69%  '

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '152.31'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.98%  '

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '5.39'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +8.19%  '

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '18.59'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +3.19%  '

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '2.76'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +19.01%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +8.46%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(44, 4).Value = '0.0₆0325'
$ws.Cells.Item(44, 5).Value = '  +15.33%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'OKB'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '42.36'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +3.15%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'USDe'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.01%  '

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '143.99'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +4.18%  '

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '3.57'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +2.81%  '

# Row 49
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '20.14'
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +6.51%  '

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '0.597'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +3.16%  '

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '0.0514'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +3.88%  '
